$wb = $excel.ActiveWorkbook

# --- Overview sheet: update "Latest HO Xliff Generate Date" for the two
#     rows that were just handed off (a0f3f310... and ff0b07d8...) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2017-02-28 07:05:47"
$overview.Range("G4").Value = "2017-02-28 07:05:47"

# --- zh-cn sheet: Priority goes from "ht" to "mt" and the Latest Handoff
#     Datetime is refreshed for the same two rows ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2017-02-28 07:05:32"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H4").Value = "2017-02-28 07:05:32"

# --- de-de sheet: Priority goes from "ht" to "mt" for the same two rows ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "mt"
$dede.Range("E4").Value = "mt"
